$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 22:52"

# Refreshed COVID-19 case numbers (columns B:H) for the affected country rows
# Row 4
$ws.Range("B4").Value = 462180
$ws.Range("C4").Value = 27253
$ws.Range("D4").Value = 24961
$ws.Range("E4").Value = 420775
$ws.Range("F4").Value = 9823
$ws.Range("G4").Value = 1656
$ws.Range("H4").Value = 16444

# Row 16
$ws.Range("B16").Value = 20703
$ws.Range("C16").Value = 1265
$ws.Range("D16").Value = 5218
$ws.Range("E16").Value = 14981
$ws.Range("F16").Value = 518
$ws.Range("G16").Value = 77
$ws.Range("H16").Value = 504

# Row 94
$ws.Range("B94").Value = 410
$ws.Range("C94").Value = 68
$ws.Range("D94").Value = 40
$ws.Range("E94").Value = 359
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 11

# Row 95
$ws.Range("B95").Value = 409
$ws.Range("C95").Value = 9
$ws.Range("D95").Value = 165
$ws.Range("E95").Value = 221
$ws.Range("F95").Value = 7
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 23

# Row 96
$ws.Range("B96").Value = 380
$ws.Range("C96").Value = 1
$ws.Range("D96").Value = 80
$ws.Range("E96").Value = 295
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 5

# Row 97
$ws.Range("B97").Value = 372
$ws.Range("C97").Value = 14
$ws.Range("D97").Value = 161
$ws.Range("E97").Value = 204
$ws.Range("F97").Value = 5
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 7

# Row 98
$ws.Range("B98").Value = 362
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 40
$ws.Range("E98").Value = 322
$ws.Range("F98").Value = 4
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0

# Row 99
$ws.Range("B99").Value = 343
$ws.Range("C99").Value = 31
$ws.Range("D99").Value = 6
$ws.Range("E99").Value = 314
$ws.Range("F99").Value = 10
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 23

# Row 163
$ws.Range("D163").Value = 2
$ws.Range("E163").Value = 18

# Row 182
$ws.Range("C182").Value = 4
$ws.Range("E182").Value = 11
$ws.Range("H182").Value = 0

# Row 183
$ws.Range("D183").Value = 1
$ws.Range("E183").Value = 10
$ws.Range("H183").Value = 1

# Row 184
$ws.Range("B184").Value = 12
$ws.Range("D184").Value = 7
$ws.Range("E184").Value = 5

# Row 186
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 11

# Row 187
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 2
$ws.Range("E187").Value = 9
$ws.Range("H187").Value = 0

# Row 188
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 8
$ws.Range("H188").Value = 3

# Row 189
$ws.Range("B189").Value = 11
$ws.Range("D189").Value = 11
$ws.Range("E189").Value = 0
$ws.Range("H189").Value = 0

# Row 190
$ws.Range("B190").Value = 10
$ws.Range("D190").Value = 4
$ws.Range("E190").Value = 5
$ws.Range("H190").Value = 1

# Row 191
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 1
$ws.Range("F191").Value = 0
$ws.Range("H191").Value = 0

# Row 192
$ws.Range("C192").Value = 1
$ws.Range("E192").Value = 8
$ws.Range("F192").Value = 1
$ws.Range("H192").Value = 1

# Row 193
$ws.Range("B193").Value = 9
$ws.Range("E193").Value = 7
$ws.Range("H193").Value = 2

# Row 194
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 8

# Row 195
$ws.Range("F195").Value = 1

# Row 196
$ws.Range("F196").Value = 0

Write-Output "Countries & provincias Spain data refreshed."
